$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "F2" = 2.84; "G2" = 3.2; "H2" = 2.48; "I2" = 2.78; "J2" = 3.25;
    "V2" = 1.56; "Y2" = 13.5; "Z2" = 21; "AA2" = 46; "AB2" = 14.5;
    "AC2" = 9.6; "AD2" = 14.5; "AE2" = 36; "AF2" = 25; "AG2" = 16;
    "AH2" = 21; "AI2" = 50; "AJ2" = 60; "AK2" = 42; "AL2" = 55;
    "AN2" = 38; "AO2" = 29;

    "I3" = 4.6; "L3" = 1.29; "N3" = 5.6; "P3" = 2.52; "Q3" = 1.62;
    "R3" = 1.61; "S3" = 2.54; "T3" = 1.63; "U3" = 2.5; "V3" = 1.27;
    "Y3" = 23; "AC3" = 9.6; "AK3" = 16.5;

    "H4" = 1.3; "J4" = 5.5; "K4" = 6.6; "L4" = 1.29; "T4" = 1.99;
    "W4" = 1.08; "X4" = 30; "AB4" = 46; "AC4" = 16.5; "AD4" = 13;
    "AF4" = 130; "AG4" = 50; "AH4" = 36; "AJ4" = 500; "AK4" = 220;
    "AL4" = 170; "AM4" = 190; "AN4" = 280; "AO4" = 6.2;

    "F5" = 1.84; "G5" = 1.96; "H5" = 4; "I5" = 4.8; "J5" = 3.95;
    "L5" = 1.27; "N5" = 4.8; "P5" = 2.32; "Q5" = 1.6; "R5" = 1.53;
    "S5" = 2.52; "U5" = 2.32; "W5" = 2.04; "AB5" = 1000; "AC5" = 10.5;
    "AD5" = 1000; "AF5" = 1000; "AG5" = 11; "AH5" = 1000; "AK5" = 1000;
    "AN5" = 9.4;

    "F6" = 1.38; "H6" = 9.4; "I6" = 9.6; "M6" = 1.03; "R6" = 1.79;
    "S6" = 2.2; "T6" = 1.77; "W6" = 3.6;

    "G7" = 3.75; "L7" = 1.32; "P7" = 2.38; "T7" = 1.62; "U7" = 2.48;
    "W7" = 1.36; "AB7" = 18; "AJ7" = 70;

    "F8" = 1.7; "G8" = 1.71; "L8" = 1.41; "N8" = 3.85; "Q8" = 1.99;
    "S8" = 3.5; "T8" = 1.97; "W8" = 2.4; "X8" = 14; "Y8" = 19.5;
    "AF8" = 9.4; "AG8" = 9.800000000000001; "AK8" = 17.5;

    "G9" = 470; "H9" = 1.01;

    "G12" = 2.14; "M12" = 1.07; "N12" = 3.55; "O12" = 1.32; "Q12" = 2.02;
    "S12" = 3.5; "T12" = 1.81; "U12" = 2.06; "W12" = 1.87; "X12" = 14.5;
    "Y12" = 15.5; "Z12" = 32; "AB12" = 9.6; "AC12" = 8.4; "AD12" = 18;
    "AE12" = 60; "AF12" = 13.5; "AG12" = 11.5; "AH12" = 19.5; "AI12" = 65;
    "AJ12" = 26; "AK12" = 24; "AL12" = 970; "AM12" = 130; "AN12" = 16.5;
    "AO12" = 65;

    "I13" = 2.8; "J13" = 3.8; "L13" = 1.3; "R13" = 1.57;

    "H14" = 5.1; "J14" = 4.5; "L14" = 1.27;
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
